# Commit: swap the deck's colour theme (the custom "Integral" / Red
# Violet design) over to the stock "Office Theme" colours, and re-style
# the B1/B2 table on slide 5 with a different (built-in) table style.

function Get-RGBValue($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table on slide 5 (the "Type of document / Definition / Why it is
#    important" table) switches from the custom "Table_0" style to the
#    built-in table style {EF296B93-1952-4EAC-92DC-19A6DCA92E15}.
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{EF296B93-1952-4EAC-92DC-19A6DCA92E15}")
}

# ---------------------------------------------------------------------
# 2) The presentation's theme colours move from the custom "Red Violet"
#    palette (theme "Integral") to the default "Office" palette (theme
#    "Office Theme"). Push the new RGB values through the theme colour
#    scheme, in theme colour order:
#      1 dk1, 2 lt1, 3 dk2, 4 lt2,
#      5-10 accent1..accent6, 11 hlink, 12 folHlink
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$scheme = $slide1.ThemeColorScheme

$scheme.Colors(1).RGB  = Get-RGBValue 0x00 0x00 0x00   # dk1
$scheme.Colors(2).RGB  = Get-RGBValue 0xFF 0xFF 0xFF   # lt1
$scheme.Colors(3).RGB  = Get-RGBValue 0x44 0x54 0x6A   # dk2
$scheme.Colors(4).RGB  = Get-RGBValue 0xE7 0xE6 0xE6   # lt2
$scheme.Colors(5).RGB  = Get-RGBValue 0x5B 0x9B 0xD5   # accent1
$scheme.Colors(6).RGB  = Get-RGBValue 0xED 0x7D 0x31   # accent2
$scheme.Colors(7).RGB  = Get-RGBValue 0xA5 0xA5 0xA5   # accent3
$scheme.Colors(8).RGB  = Get-RGBValue 0xFF 0xC0 0x00   # accent4
$scheme.Colors(9).RGB  = Get-RGBValue 0x44 0x72 0xC4   # accent5
$scheme.Colors(10).RGB = Get-RGBValue 0x70 0xAD 0x47   # accent6
$scheme.Colors(11).RGB = Get-RGBValue 0x05 0x63 0xC1   # hlink
$scheme.Colors(12).RGB = Get-RGBValue 0x95 0x4F 0x72   # folHlink
